$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.941.30'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").Value = '1.907.55'
$ws.Range("E3").Value = '  +0.78%  '
$ws.Range("D4").Value = '''0.9988'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '''0.8044'
$ws.Range("E5").Value = '  +6.05%  '
$ws.Range("D6").Value = '''241.74'
$ws.Range("E6").Value = '  +1.03%  '
$ws.Range("D7").Value = '''0.9993'
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").Value = '''0.3157'
$ws.Range("E8").Value = '  +3.84%  '
$ws.Range("D9").Value = '''26.38'
$ws.Range("E9").Value = '  +4.69%  '
$ws.Range("D10").Value = '''0.06907'
$ws.Range("E10").Value = '  +1.35%  '
$ws.Range("D11").Value = '''0.07988'
$ws.Range("E11").Value = '  +0.05%  '
$ws.Range("D12").Value = '1.911.27'
$ws.Range("E12").Value = '  +1.05%  '
$ws.Range("D13").Value = '''0.7384'
$ws.Range("E13").Value = '  -1.20%  '
$ws.Range("D14").Value = '''5.190'
$ws.Range("E14").Value = '  -0.08%  '
$ws.Range("E15").Value = '  +2.02%  '
$ws.Range("D16").Value = '29.934.38'
$ws.Range("E16").Value = '  +0.39%  '
$ws.Range("D17").Value = '''13.97'
$ws.Range("E17").Value = '  +0.94%  '
$ws.Range("D18").Value = '''5.870'
$ws.Range("E18").Value = '  -1.83%  '
$ws.Range("D19").Value = '''245.17'
$ws.Range("E19").Value = '  +4.30%  '
$ws.Range("D20").Value = '''0.000007735'
$ws.Range("E20").Value = '  +1.02%  '
$ws.Range("D21").Value = '''0.9990'
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("D22").Value = '2.148.78'
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("D23").Value = '''0.9984'
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D24").Value = '''6.821'
$ws.Range("E24").Value = '  -1.37%  '
$ws.Range("D25").Value = '''167.69'
$ws.Range("E25").Value = '  +1.47%  '
$ws.Range("D26").Value = '''9.197'
$ws.Range("E26").Value = '  -0.45%  '
$ws.Range("D27").Value = '''0.1415'
$ws.Range("E27").Value = '  +10.22%  '
$ws.Range("E28").Value = '  +1.15%  '
$ws.Range("D29").Value = '''2.032'
$ws.Range("E29").Value = '  -0.97%  '
$ws.Range("D30").Value = '''1.363'
$ws.Range("E30").Value = '  +1.66%  '
$ws.Range("D31").Value = '''1.514'
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("D32").Value = '''4.300'
$ws.Range("E32").Value = '  +0.73%  '
$ws.Range("D33").Value = '''4.082'
$ws.Range("E33").Value = '  +1.69%  '
$ws.Range("D34").Value = '''0.05474'
$ws.Range("E34").Value = '  +2.34%  '
$ws.Range("D35").Value = '''1.264'
$ws.Range("E35").Value = '  +1.70%  '
$ws.Range("D36").Value = '''0.7293'
$ws.Range("E36").Value = '  +0.38%  '
$ws.Range("D37").Value = '''2.716'
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("E38").Value = '  -0.32%  '
$ws.Range("D39").Value = '''2.778'
$ws.Range("E39").Value = '  +0.38%  '
$ws.Range("D40").Value = '''6.134'
$ws.Range("E40").Value = '  -0.82%  '
$ws.Range("D41").Value = '''0.4416'
$ws.Range("E41").Value = '  +0.38%  '
$ws.Range("D42").Value = '''72.26'
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("D43").Value = '''0.9991'
$ws.Range("E43").Value = '  -0.17%  '
$ws.Range("D44").Value = '''0.8333'
$ws.Range("E44").Value = '  +1.33%  '
$ws.Range("D45").Value = '''1.875'
$ws.Range("E45").Value = '  -2.06%  '
$ws.Range("D46").Value = '''100.39'
$ws.Range("E46").Value = '  -0.66%  '
$ws.Range("D47").Value = '''7.533'
$ws.Range("E47").Value = '  -0.56%  '
$ws.Range("D48").Value = '''9.769'
$ws.Range("E48").Value = '  -0.81%  '
$ws.Range("D49").Value = '''985.49'
$ws.Range("E49").Value = '  +7.98%  '
$ws.Range("D50").Value = '2.055.60'
$ws.Range("E50").Value = '  +0.59%  '
$ws.Range("D51").Value = '''36.21'
$ws.Range("E51").Value = '  +0.43%  '
